$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4887.095
$ws.Range("I40").Value = 2624
$ws.Range("K40").Value = 2624
$ws.Range("M40").Value = -2449
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
$ws.Range("H132").Value = 1480.28
$ws.Range("I132").Value = 1396.125
$ws.Range("K132").Value = 4188.375
$ws.Range("M132").Value = -1658.375
$ws.Range("H138").Value = 4134.913
$ws.Range("I138").Value = 2535.889
$ws.Range("J138").Value = 5162.857
$ws.Range("K138").Value = 7607.667
$ws.Range("L138").Value = 15488.571
$ws.Range("M138").Value = -2467.667
$ws.Range("N138").Value = -25768.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7028.516
$ws.Range("I61").Value = 7534.1924
$ws.Range("K61").Value = 7534.1924
$ws.Range("M61").Value = -7322.1924
$ws.Range("H102").Value = 2998.5652
$ws.Range("I102").Value = 1939.4706
$ws.Range("J102").Value = 5999.3335
$ws.Range("K102").Value = 1939.4706
$ws.Range("L102").Value = 5999.3335
$ws.Range("M102").Value = -317.4706000000001
$ws.Range("N102").Value = -9243.333500000001
$ws.Range("H136").Value = 7028.516
$ws.Range("I136").Value = 7534.1924
$ws.Range("K136").Value = 22602.5772
$ws.Range("M136").Value = -20052.5772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69900
$ws.Range("J2").Value = 69900
$ws.Range("L2").Value = 69900
$ws.Range("N2").Value = -70126
$ws.Range("H20").Value = 5107.923
$ws.Range("I20").Value = 5698.2
$ws.Range("K20").Value = 5698.2
$ws.Range("M20").Value = -5451.2
$ws.Range("H26").Value = 20235.5
$ws.Range("I26").Value = 20235.5
$ws.Range("K26").Value = 20235.5
$ws.Range("M26").Value = -19943.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 969.25
$ws.Range("I94").Value = 829.35
$ws.Range("K94").Value = 829.35
$ws.Range("M94").Value = -378.35
$ws.Range("H96").Value = 11853.25
$ws.Range("I96").Value = 11853.25
$ws.Range("K96").Value = 11853.25
$ws.Range("M96").Value = -9107.25
$ws.Range("H105").Value = 1004.375
$ws.Range("I105").Value = 871.4
$ws.Range("K105").Value = 871.4
$ws.Range("M105").Value = 875.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 3637.1667
$ws.Range("I11").Value = 6833.3335
$ws.Range("K11").Value = 6833.3335
$ws.Range("M11").Value = -6693.3335
$ws.Range("H31").Value = 6725.2856
$ws.Range("I31").Value = 6698
$ws.Range("J31").Value = 6733.8125
$ws.Range("K31").Value = 6698
$ws.Range("L31").Value = 6733.8125
$ws.Range("M31").Value = -6403
$ws.Range("N31").Value = -7323.8125
$ws.Range("H34").Value = 6725.2856
$ws.Range("I34").Value = 6698
$ws.Range("J34").Value = 6733.8125
$ws.Range("K34").Value = 6698
$ws.Range("L34").Value = 6733.8125
$ws.Range("M34").Value = -6496
$ws.Range("N34").Value = -7137.8125
$ws.Range("H58").Value = 12601.625
$ws.Range("I58").Value = 1399
$ws.Range("J58").Value = 14202
$ws.Range("K58").Value = 1399
$ws.Range("L58").Value = 14202
$ws.Range("M58").Value = -1196
$ws.Range("N58").Value = -14608
$ws.Range("H105").Value = 1159.8889
$ws.Range("I105").Value = 1198.625
$ws.Range("J105").Value = 850
$ws.Range("K105").Value = 1198.625
$ws.Range("L105").Value = 850
$ws.Range("M105").Value = 548.375
$ws.Range("N105").Value = -4344
$ws.Range("H136").Value = 12601.625
$ws.Range("I136").Value = 1399
$ws.Range("J136").Value = 14202
$ws.Range("K136").Value = 4197
$ws.Range("L136").Value = 42606
$ws.Range("M136").Value = -1647
$ws.Range("N136").Value = -47706
$ws.Range("H139").Value = 89499.5
$ws.Range("J139").Value = 18999
$ws.Range("L139").Value = 18999
$ws.Range("N139").Value = -29279
$ws.Range("H141").Value = 213099.7
$ws.Range("J141").Value = 213099.7
$ws.Range("L141").Value = 213099.7
$ws.Range("N141").Value = -223459.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 709.9
$ws.Range("I86").Value = 2200
$ws.Range("J86").Value = 337.375
$ws.Range("K86").Value = 6600
$ws.Range("L86").Value = 1012.125
$ws.Range("M86").Value = -5414
$ws.Range("N86").Value = -3384.125
$ws.Range("H89").Value = 709.9
$ws.Range("I89").Value = 2200
$ws.Range("J89").Value = 337.375
$ws.Range("K89").Value = 19800
$ws.Range("L89").Value = 3036.375
$ws.Range("M89").Value = -13872
$ws.Range("N89").Value = -14892.375
$ws.Range("H131").Value = 12822825
$ws.Range("J131").Value = 2644.4375
$ws.Range("L131").Value = 7933.3125
$ws.Range("N131").Value = -18013.3125
$ws.Range("H140").Value = 1627.1875
$ws.Range("I140").Value = 1323.9286
$ws.Range("K140").Value = 3971.7858
$ws.Range("M140").Value = 1208.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 17993.8
$ws.Range("J70").Value = 17993.8
$ws.Range("L70").Value = 17993.8
$ws.Range("N70").Value = -18533.8
$ws.Range("H73").Value = 17993.8
$ws.Range("J73").Value = 17993.8
$ws.Range("L73").Value = 17993.8
$ws.Range("N73").Value = -19865.8
$ws.Range("H80").Value = 2499.1428
$ws.Range("J80").Value = 2666
$ws.Range("L80").Value = 2666
$ws.Range("N80").Value = -4662
$ws.Range("H83").Value = 2499.1428
$ws.Range("J83").Value = 2666
$ws.Range("L83").Value = 13330
$ws.Range("N83").Value = -23314
$ws.Range("H97").Value = 882.94116
$ws.Range("I97").Value = 743.5
$ws.Range("K97").Value = 743.5
$ws.Range("M97").Value = -247.5
$ws.Range("H102").Value = 5031.4814
$ws.Range("I102").Value = 2833.0557
$ws.Range("K102").Value = 2833.0557
$ws.Range("M102").Value = -1211.0557
$ws.Range("H132").Value = 4403.8076
$ws.Range("I132").Value = 3975.6191
$ws.Range("J132").Value = 6202.2
$ws.Range("K132").Value = 11926.8573
$ws.Range("L132").Value = 18606.6
$ws.Range("M132").Value = -9396.8573
$ws.Range("N132").Value = -23666.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 21420.154
$ws.Range("I93").Value = 6223.4546
$ws.Range("K93").Value = 6223.4546
$ws.Range("M93").Value = -4975.4546
$ws.Range("H122").Value = 2177.182
$ws.Range("I122").Value = 2177.182
$ws.Range("K122").Value = 6531.545999999999
$ws.Range("M122").Value = -4081.545999999999
$ws.Range("H132").Value = 15379.333
$ws.Range("I132").Value = 16049.286
$ws.Range("K132").Value = 48147.858
$ws.Range("M132").Value = -45617.858
$ws.Range("H136").Value = 5424.926
$ws.Range("I136").Value = 5424.926
$ws.Range("K136").Value = 16274.778
$ws.Range("M136").Value = -13724.778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7250.25
$ws.Range("J62").Value = 9999.5
$ws.Range("L62").Value = 9999.5
$ws.Range("N62").Value = -11247.5
$ws.Range("H65").Value = 7250.25
$ws.Range("J65").Value = 9999.5
$ws.Range("L65").Value = 49997.5
$ws.Range("N65").Value = -56237.5
$ws.Range("H81").Value = 2086.6
$ws.Range("I81").Value = 2086.6
$ws.Range("K81").Value = 4173.2
$ws.Range("M81").Value = -3112.2
$ws.Range("H84").Value = 2086.6
$ws.Range("I84").Value = 2086.6
$ws.Range("K84").Value = 20866
$ws.Range("M84").Value = -15562
$ws.Range("H107").Value = 1713
$ws.Range("I107").Value = 965.2
$ws.Range("J107").Value = 2180.375
$ws.Range("K107").Value = 2895.6
$ws.Range("L107").Value = 6541.125
$ws.Range("M107").Value = -975.6000000000004
$ws.Range("N107").Value = -10381.125
$ws.Range("H122").Value = 2637.8333
$ws.Range("I122").Value = 2055.0334
$ws.Range("J122").Value = 4094.8333
$ws.Range("K122").Value = 6165.100199999999
$ws.Range("L122").Value = 12284.4999
$ws.Range("M122").Value = -3715.100199999999
$ws.Range("N122").Value = -17184.4999
$ws.Range("H132").Value = 3987.309
$ws.Range("I132").Value = 3339.2954
$ws.Range("J132").Value = 6579.364
$ws.Range("K132").Value = 10017.8862
$ws.Range("L132").Value = 19738.092
$ws.Range("M132").Value = -7487.886200000001
$ws.Range("N132").Value = -24798.092
$ws.Range("H136").Value = 3118.4856
$ws.Range("I136").Value = 2876.75
$ws.Range("K136").Value = 8630.25
$ws.Range("M136").Value = -6080.25
